# Insert a new data row at row 239 (shifts existing rows 239-297 down to 240-298)
# and populate it with a new weekly price record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(239).Insert()

$ws.Range("A239").Value = 8
$ws.Range("B239").Value = "Terminal La Palmera de La Serena"
$ws.Range("C239").Value = "Coquimbo"
$ws.Range("D239").Value = 44889
$ws.Range("E239").Value = 4
$ws.Range("F239").Value = 100112031
$ws.Range("G239").Value = "Poroto verde"
$ws.Range("H239").Value = "Magnum"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 440
$ws.Range("K239").Value = 31000
$ws.Range("L239").Value = 32000
$ws.Range("M239").Value = 31500
$ws.Range("N239").Value = "`$/malla 25 kilos"
$ws.Range("O239").Value = "Provincia de Limarí"
$ws.Range("P239").Value = 1260
$ws.Range("Q239").Value = 25
$ws.Range("R239").Value = "Hortaliza"
